# Update the "Raju Ahamed" cash sheet for 19.06.19 (today's sales update).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raju Ahamed")

# Update the date header shown in both copies of the sheet (top A4:J4 and
# bottom A31:J31 merged ranges share the same underlying text).
$ws.Range("A4").Value = "Date: 19.06.19"
$ws.Range("A31").Value = "Date: 19.06.19"

# --- Top table (rows 6-11): update today's quantities sold ---
$ws.Range("E6").Value = 22
$ws.Range("E7").Value = 113
$ws.Range("E8").Value = 303
$ws.Range("E9").Value = 16
$ws.Range("E10").Value = 20

# --- Bottom duplicate table (rows 33-38): same quantities ---
$ws.Range("E33").Value = 22
$ws.Range("E34").Value = 113
$ws.Range("E35").Value = 303
$ws.Range("E36").Value = 16
$ws.Range("E37").Value = 20
